# Simplify database schema migration
# Appends one new data row (row 76) to each of the four sheets
# (FE_LFT_#1, FE_LFT_#2, FE_PLT_#1, FE_PLT_#2), mirroring the existing
# row-75 layout/format and updating the used dimension accordingly.

$wb = $excel.ActiveWorkbook

$rowsData = @{
    1 = @{
        A = 45862.49133101852
        B = "0x01,0x7c"
        C = "0x00,0xa6,0x46,0x93,0x3c,0x23,0x3f,0x43,0xe8,0xa0,"
        D = "0x01,0x28"
        E = "0xf"
        F = 380
        G = [double]"7.598631275147109e+23"
        H = 296
        I = 15
    }
    2 = @{
        A = 45862.49133101852
        B = "0x01,0x90"
        C = "0x00,0xa6,0x60,0x33,0x96,0x39,0x62,0xd0,0x5e,0x78,"
        D = "0x01,0x38"
        E = "0xe"
        F = 400
        G = [double]"5.68432987514711e+23"
        H = 312
        I = 14
    }
    3 = @{
        A = 45862.49133101852
        B = "0x00,0x6e"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"
        D = "0x00,0x61"
        E = "0x3"
        F = 110
        G = [double]"5.68631262647114e+23"
        H = 97
        I = 3
    }
    4 = @{
        A = 45862.49133101852
        B = "0x00,0x6e"
        C = "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c,"
        D = "0x00,0x60"
        E = "0x3"
        F = 110
        G = [double]"9.85046333984776e+23"
        H = 96
        I = 3
    }
}

foreach ($sheetIdx in 1..4) {
    $ws = $wb.Worksheets.Item($sheetIdx)
    $data = $rowsData[$sheetIdx]

    # Duplicate row 75's formatting/layout into the new row 76, then
    # overwrite with the new row's own values.
    $ws.Range("A75:I75").Copy($ws.Range("A76:I76"))

    $ws.Cells.Item(76, 1).Value = $data.A
    $ws.Cells.Item(76, 2).Value = $data.B
    $ws.Cells.Item(76, 3).Value = $data.C
    $ws.Cells.Item(76, 4).Value = $data.D
    $ws.Cells.Item(76, 5).Value = $data.E
    $ws.Cells.Item(76, 6).Value = $data.F
    $ws.Cells.Item(76, 7).Value = $data.G
    $ws.Cells.Item(76, 8).Value = $data.H
    $ws.Cells.Item(76, 9).Value = $data.I
}

"Appended row 76 to all sheets"
